$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Relocate existing cells that keep their original text, capturing style
#    (Copy(Destination) duplicates both value and style/format)
# ---------------------------------------------------------------------------
# old D1 "Raw Code (picture)" (s1)  -> F1
$ws.Range("D1").Copy($ws.Range("F1"))
# old D2 hyperlinked image url (s2) -> F2
$ws.Range("D2").Copy($ws.Range("F2"))
# old C2 "Epoch 10/10..." (s4)      -> D2 (same letter, new row position meaning)
$ws.Range("C2").Copy($ws.Range("D2"))
# old C2 "Epoch 10/10..." (s4)      -> D3 (second copy; text replaced later)
$ws.Range("C2").Copy($ws.Range("D3"))

# ---------------------------------------------------------------------------
# 2) Grab formatting (style) for brand-new header cells from old C1, before
#    C1's own text gets overwritten.
# ---------------------------------------------------------------------------
$ws.Range("C1").Copy($ws.Range("D1"))
$ws.Range("C1").Copy($ws.Range("E1"))
# old B1 "Model Description" (s5)   -> C1 (new Preformance header)
$ws.Range("B1").Copy($ws.Range("C1"))

# ---------------------------------------------------------------------------
# 3) Grab formatting for the new row 3 description cells from old B2, before
#    B2's own text gets overwritten.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy($ws.Range("B3"))
$ws.Range("B2").Copy($ws.Range("C3"))

# Old C2 is no longer needed on its own (its content now lives at D2/D3).
$ws.Range("C2").Clear()

# ---------------------------------------------------------------------------
# 4) Create the combined "hyperlink + wrap" style (used by E2/E3) by copying
#    the existing hyperlink-styled cell and then enabling wrap text.
# ---------------------------------------------------------------------------
$ws.Range("D2").Copy($ws.Range("E2"))
$ws.Range("E2").WrapText = $true
$ws.Range("E2").Copy($ws.Range("E3"))

# ---------------------------------------------------------------------------
# 5) Write the new text content, in an order chosen to reproduce the target
#    shared-string table ordering.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Training Details"
$ws.Range("E1").Value = "Summary"
$ws.Range("E2").Value = "https://i.imgur.com/VQoLgxR.png"
$ws.Range("A3").Value = "models/04_20_3:59:37PM/"
$ws.Range("B2").Value = "The default model from the tutoral"
$ws.Range("B3").Value = "Default model with layser.Dropout(0.2)"
$ws.Range("D3").Value = "loss: 0.0508 - accuracy: 0.9812 - val_loss: 1.2555 - val_accuracy: 0.8237"
$ws.Range("E3").Value = "https://i.imgur.com/EYitGMG.png"
$ws.Range("C1").Value = "Preformance"
$ws.Range("C3").Value = "Worse than row two predicting everything with 11% certanty"

# ---------------------------------------------------------------------------
# 6) Hyperlinks (values already match display text, so no TextToDisplay
#    override is required).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://i.imgur.com/xPgij9j.png")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://i.imgur.com/VQoLgxR.png")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://i.imgur.com/EYitGMG.png")

# ---------------------------------------------------------------------------
# 7) Column widths (runtime quantizes ColumnWidth to 1/6-character pixel
#    steps, so these inputs are chosen to land as close as possible to the
#    target widths).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.333333333333332
$ws.Columns.Item(2).ColumnWidth = 26.5
$ws.Columns.Item(3).ColumnWidth = 26.5
$ws.Columns.Item(4).ColumnWidth = 21.833333333333332
$ws.Columns.Item(5).ColumnWidth = 14
$ws.Columns.Item(6).ColumnWidth = 28.166666666666668

# ---------------------------------------------------------------------------
# 8) Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 85
$ws.Rows.Item(3).RowHeight = 51

# ---------------------------------------------------------------------------
# 9) Sheet view settings
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.DisplayFormulas = $true
$win.Zoom = 109
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("E3").Select()
